$wb = $excel.ActiveWorkbook

# --- Rename the original sheet to "My pc" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "My pc"

# --- Duplicate it (placed right after "My pc") and rename the copy to "ComputeCanada" ---
$null = $ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "ComputeCanada"

# "ComputeCanada" only keeps the first data block (rows 1-7); drop the extra
# pre-formatted blank rows that were copied from "My pc" (rows 8-26).
$ws2.Range("A8:J26").EntireRow.Delete()

# This tab hasn't been run yet, so its results columns are still empty.
$ws2.Range("G3:J7").ClearContents()

# --- Column widths (slightly widened on both tabs) ---
$ws1.Columns("A:F").ColumnWidth = 10.6
$ws1.Columns("G:J").ColumnWidth = 11.3
$ws2.Columns("A:F").ColumnWidth = 10.6
$ws2.Columns("G:J").ColumnWidth = 11.3

# --- Fill in the newly measured "My pc" results for row 7 ---
$ws1.Range("I7").Value = 170597
$ws1.Range("J7").Value = 46204

# --- Selections on each tab ---
$null = $ws1.Range("G11").Select()
$null = $ws2.Range("I11").Select()

# --- ComputeCanada is the active tab ---
$null = $ws2.Activate()
